# Automatische test-sync: 2025-06-22 21:46:50
# Append a new log row (row 46) to the "Logs" sheet and bump the
# "Sollicitatie / Vacature" count on the "Dashboard" sheet.

$wb = $excel.ActiveWorkbook

$logs = $wb.Worksheets.Item("Logs")

$newRow = 46

$logs.Cells.Item($newRow, 1).Value = "Sollicitatie salesfunctie"
$logs.Cells.Item($newRow, 2).Value = "mailmind.test@zohomail.eu"
$logs.Cells.Item($newRow, 3).Value = "Hierbij mijn sollicitatie voor de salesfunctie. CV in bijlage."
$logs.Cells.Item($newRow, 4).Value = "Sollicitatie / Vacature"
$logs.Cells.Item($newRow, 5).Value = "Beste sollicitant,`nDank u voor uw interesse in de salesfunctie bij ons bedrijf. We hebben uw sollicitatie en CV ontvangen. Onze HR-afdeling zal uw sollicitatie zo spoedig mogelijk bekijken en contact met u opnemen indien uw profiel aansluit bij de functievereisten.`nMet vriendelijke groet,`n[Naam bedrijf]"
$logs.Cells.Item($newRow, 6).Value = "2025-06-22 21:45:50"
$logs.Cells.Item($newRow, 7).Value = "Ja"

# The multi-line Antwoord text makes the runtime auto-size the row; put it
# back to the sheet default (and drop the explicit/custom height) so the
# row matches the rest of the sheet.
$logs.Rows.Item($newRow).EntireRow.AutoFit()

# Extend conditional formatting ranges to include the new row, in place
# (keeps rule type/operator/dxf/priority untouched, only grows the sqref).
$dRange = $logs.Range("D2:D45")
for ($i = 1; $i -le $dRange.FormatConditions.Count; $i++) {
    $dRange.FormatConditions.Item($i).ModifyAppliesToRange($logs.Range("D2:D46"))
}

$gRange = $logs.Range("G2:G45")
for ($i = 1; $i -le $gRange.FormatConditions.Count; $i++) {
    $gRange.FormatConditions.Item($i).ModifyAppliesToRange($logs.Range("G2:G46"))
}

# Update the Dashboard count for "Sollicitatie / Vacature".
$dashboard = $wb.Worksheets.Item("Dashboard")
$dashboard.Cells.Item(3, 2).Value = 6
